$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Rows 22-27: the date in column C moves from 15-May-2020 (43966) to 18-May-2020 (43969)
$ws.Range("C22").Value = 43969
$ws.Range("C23").Value = 43969
$ws.Range("C24").Value = 43969
$ws.Range("C25").Value = 43969
$ws.Range("C26").Value = 43969
$ws.Range("C27").Value = 43969

# New entries for 19-May-2020 (43970) filling in previously-empty rows 28-31
$ws.Range("A28").Value = "SD02+DCD02"
$ws.Range("C28").Value = 43970
$ws.Range("D28").Value = 0.354166666666667
$ws.Range("E28").Value = 0.416666666666667

$ws.Range("A29").Value = "Lav UC03 og review UC10"
$ws.Range("C29").Value = 43970
$ws.Range("D29").Value = 0.416666666666667
$ws.Range("E29").Value = 0.541666666666667

$ws.Range("A30").Value = "Vejledning fra Anders"
$ws.Range("C30").Value = 43970
$ws.Range("D30").Value = 0.541666666666667
$ws.Range("E30").Value = 0.59375

$ws.Range("A31").Value = "SD01xx-DCD01xx"
$ws.Range("C31").Value = 43970
$ws.Range("D31").Value = 0.59375
$ws.Range("E31").Value = 0.697916666666667

# Move the active selection on Ark1 to D48, matching the saved view state
$ws.Range("D48").Select()
